$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.88720612145694411
$ws.Range("BP1").Value = 0.90855576171814478
$ws.Range("A3").Value = 0.62508896000868419
$ws.Range("B3").Value = 0.7245978608935012
$ws.Range("D3").Value = 0.91952500895389266
$ws.Range("I3").Value = 0.96231827651185209
$ws.Range("AH3").Value = 0.91944729511803691
$ws.Range("BA4").Value = 0.65171846872070827
$ws.Range("BB4").Value = 0.76171155902223342
$ws.Range("F5").Value = 0.7781653324565132
$ws.Range("G5").Value = 0.67128876370545576
$ws.Range("D6").Value = 0.99729643016830405
$ws.Range("AN6").Value = 0.83607622892135969
$ws.Range("F7").Value = 0.89361643769866417
$ws.Range("H7").Value = 0.63445865593810213
$ws.Range("Q7").Value = 0.68443663198071214
$ws.Range("I8").Value = 0.99728152884249721
$ws.Range("B10").Value = 0.97709529177196719
$ws.Range("H10").Value = 0.80787963004347474
$ws.Range("I10").Value = 0.91123217021570058
$ws.Range("AW10").Value = 0.69354284898409468
$ws.Range("BH11").Value = 0.74371088705433663
$ws.Range("J12").Value = 0.78390678812140735
$ws.Range("N12").Value = 0.92950728725993115
$ws.Range("L13").Value = 0.91996168909472265
$ws.Range("M14").Value = 0.73509352543541884
$ws.Range("P14").Value = 0.99267171505010254
$ws.Range("M15").Value = 0.99444118810692639
$ws.Range("N15").Value = 0.87147427915351272
$ws.Range("R16").Value = 0.63849471901932597
$ws.Range("C17").Value = 0.95992422866566895
$ws.Range("P17").Value = 0.82961485631210907
$ws.Range("R17").Value = 0.7377675457304258
$ws.Range("AK17").Value = 0.97522260004205918
$ws.Range("K18").Value = 0.9814944691439117
$ws.Range("S18").Value = 0.80347310826140372
$ws.Range("T18").Value = 0.68478108302925411
$ws.Range("T19").Value = 0.66558519951133688
$ws.Range("U19").Value = 0.82826717140416073
$ws.Range("V20").Value = 0.89877996672932592
$ws.Range("T21").Value = 0.92796638848003044
$ws.Range("V21").Value = 0.65744508825418846
$ws.Range("W22").Value = 0.9974800561047541
$ws.Range("X22").Value = 0.6243972702265207
$ws.Range("BB22").Value = 0.75661370964110042
$ws.Range("U23").Value = 0.92023150764441919
$ws.Range("W24").Value = 0.67222041893619089
$ws.Range("Y24").Value = 0.92694354555131653
$ws.Range("AQ24").Value = 0.62687683503649494
$ws.Range("K25").Value = 0.65166079133631061
$ws.Range("W25").Value = 0.98037270557884937
$ws.Range("AA25").Value = 0.92281589006080189
$ws.Range("AR26").Value = 0.9614070074926202
$ws.Range("AB27").Value = 0.9603509349570658
$ws.Range("AX27").Value = 0.67130209428121557
$ws.Range("Z28").Value = 0.62839826263490717
$ws.Range("BI29").Value = 0.94986371661241009
$ws.Range("BO29").Value = 0.96966848331276312
$ws.Range("S30").Value = 0.97689988284417695
$ws.Range("AC30").Value = 0.97570239989733332
$ws.Range("AC31").Value = 0.55696400294780457
$ws.Range("AD31").Value = 0.88522408995933388
$ws.Range("AF31").Value = 0.91244919059023255
$ws.Range("AG32").Value = 0.79742669519091769
$ws.Range("BD32").Value = 0.87720348013935878
$ws.Range("AE33").Value = 0.96673130056794698
$ws.Range("AI33").Value = 0.73564363744498906
$ws.Range("AF34").Value = 0.64976975569040674
$ws.Range("AG34").Value = 0.96112080175375125
$ws.Range("AI34").Value = 0.7916163680988304
$ws.Range("BJ34").Value = 0.75741514438178836
$ws.Range("AJ35").Value = 0.72086906539711249
$ws.Range("AL36").Value = 0.78758746516196465
$ws.Range("AI37").Value = 0.97712184924613443
$ws.Range("AJ37").Value = 0.77410204868858934
$ws.Range("AL37").Value = 0.9019865265164142
$ws.Range("AM37").Value = 0.9796716894642713
$ws.Range("AO38").Value = 0.69593126498163271
$ws.Range("AN39").Value = 0.81190119112850834
$ws.Range("AL40").Value = 0.82084172594066573
$ws.Range("AP40").Value = 0.753945933472424
$ws.Range("AM41").Value = 0.7207852465593958
$ws.Range("AN41").Value = 0.77322781194772316
$ws.Range("AO42").Value = 0.99836140534316897
$ws.Range("AP43").Value = 0.83014994239232109
$ws.Range("AR43").Value = 0.97354825413689683
$ws.Range("AX43").Value = 0.65643281931736874
$ws.Range("AP44").Value = 0.85893691318887577
$ws.Range("AR45").Value = 0.85619166117736245
$ws.Range("AR46").Value = 0.90423711951960495
$ws.Range("AS46").Value = 0.6952982975850639
$ws.Range("AU46").Value = 0.91472829715841364
$ws.Range("AV46").Value = 0.94466999097922166
$ws.Range("AA47").Value = 0.8879728148442263
$ws.Range("AU48").Value = 0.94983923328788777
$ws.Range("AW48").Value = 0.85967794896043315
$ws.Range("AX48").Value = 0.93966779795300015
$ws.Range("AS49").Value = 0.61073485558430651
$ws.Range("AY49").Value = 0.69927594045939234
$ws.Range("AW50").Value = 0.7970377748961075
$ws.Range("I51").Value = 0.96829903530705175
$ws.Range("AJ51").Value = 0.89871681405433623
$ws.Range("AX51").Value = 0.9272835629733347
$ws.Range("BB52").Value = 0.92555129159302907
$ws.Range("I53").Value = 0.63161165623308035
$ws.Range("BD54").Value = 0.77012649698010571
$ws.Range("BD55").Value = 0.54161667064742169
$ws.Range("BE55").Value = 0.92039483299138758
$ws.Range("BF56").Value = 0.73276279052578408
$ws.Range("D57").Value = 0.93528898470161403
$ws.Range("AO58").Value = 0.75954497501352025
$ws.Range("N59").Value = 0.91819718985871734
$ws.Range("AB59").Value = 0.69942359176632762
$ws.Range("BE59").Value = 0.54048973188260629
$ws.Range("BF59").Value = 0.98922716050952431
$ws.Range("BF60").Value = 0.81209470437280795
$ws.Range("BG60").Value = 0.92721391967483568
$ws.Range("BI60").Value = 0.62831368755379258
$ws.Range("AB62").Value = 0.92585161510664693
$ws.Range("BH62").Value = 0.99363614420949931
$ws.Range("BI62").Value = 0.95422675481271457
$ws.Range("BK62").Value = 0.92734346085527697
$ws.Range("BN63").Value = 0.860225573858165
$ws.Range("G64").Value = 0.96992387282362236
$ws.Range("O64").Value = 0.95218792086850046
$ws.Range("AZ64").Value = 0.6671731819735256
$ws.Range("BK64").Value = 0.78397751570876095
$ws.Range("BL66").Value = 0.97629414798286351
$ws.Range("BM66").Value = 0.89766190166603854
$ws.Range("BP66").Value = 0.71789498358167836
$ws.Range("O67").Value = 0.86099178021556844
$ws.Range("BM67").Value = 0.75364066355407844
$ws.Range("BN67").Value = 0.88857931311569471
$ws.Range("O68").Value = 0.98190104062709671
$ws.Range("AA68").Value = 0.95864583233607492
$ws.Range("AC68").Value = 0.94308075216362075
